$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the policy number text value (E2) with the new value for Irina's PC env.
# Leading apostrophe keeps it entered as text (preserves the existing
# text/quote-prefix cell format instead of Excel reinterpreting it as a number).
$ws.Range("E2").Value = "'12112001753"

# Move/update the active selection as left by the author before saving
$ws.Activate()
$ws.Range("E3").Select()
